$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 36 (pushes the existing rows 36-61 down to 37-62)
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row with the new "Membrillo / Champion" price record
$ws.Cells.Item(36, 1).Value = 6
$ws.Cells.Item(36, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(36, 3).Value = "Metropolitana"
$ws.Cells.Item(36, 4).Value = 44634
$ws.Cells.Item(36, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 5).Value = 13
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100104
$ws.Cells.Item(36, 8).Value = "Frutos de pepita"
$ws.Cells.Item(36, 9).Value = 100104003
$ws.Cells.Item(36, 10).Value = "Membrillo"
$ws.Cells.Item(36, 11).Value = "Champion"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 8
$ws.Cells.Item(36, 14).Value = 250000
$ws.Cells.Item(36, 15).Value = 250000
$ws.Cells.Item(36, 16).Value = 250000
$ws.Cells.Item(36, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(36, 18).Value = "Región Metropolitana"
$ws.Cells.Item(36, 19).Value = 556
$ws.Cells.Item(36, 20).Value = 450
